$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Birat BF")
$ws.Activate()

# Add new scenario rows (A11:A18) with base-case scenario names for several
# regions/countries, extending the "medium-complexity" steel BF model.
$scenarios = @(
    "global-BF-base",
    "China-BF-base",
    "EU-BF-base",
    "India-BF-base",
    "Japan-BF-base",
    "Russia-BF-base",
    "SouthKorea-BF-base",
    "USA-BF-base"
)

$row = 11
foreach ($name in $scenarios) {
    $ws.Cells.Item($row, 1).Value = $name
    $row = $row + 1
}

# Update the selected range on the active pane to reflect the new data entry area
$ws.Range("A1:A4").Select()
